$d = $word.ActiveDocument

$newText = "Kampagnendaten Pegasus: 8. bis 17. Oktober, 7. bis 16. November,"

# Collect the paragraphs that hold the old "Kampagnendaten 2018 ... Perseus ..."
# text (there are four occurrences in this document).
$targets = New-Object System.Collections.ArrayList
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Kampagnendaten*") {
        $targets.Add($p) | Out-Null
    }
}

# Walk backwards so earlier character offsets stay valid while we edit.
for ($i = $targets.Count - 1; $i -ge 0; $i--) {
    $p = $targets[$i]
    $start = $p.Range.Start
    $end = $p.Range.End - 1   # exclude the trailing paragraph mark

    # Clear every run in the paragraph (this removes the runs entirely,
    # rather than leaving one run behind with inherited direct formatting).
    $clearRange = $d.Range($start, $end)
    $clearRange.Text = ""

    # Insert the new text into the now-empty paragraph; a run created in an
    # empty paragraph like this carries no direct run formatting (no rPr),
    # matching the target markup of a single plain <w:r><w:t>...</w:t></w:r>.
    $insertRange = $d.Range($start, $start)
    $insertRange.InsertAfter($newText)
}
